$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roberto")

# Copy the formatting (style) of A1 down onto the new label cells A4:A9
$ws.Range("A1").Copy()
$ws.Range("A4:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$data = @(
    @("Puntuación general airbnb de Clara: ", 5),
    @("Cantidad reseñas airbnb de Clara: ", 204),
    @("Precio por noche airbnb de Clara: ", 46),
    @("Puntuación general airbnb de Roberto: ", 4),
    @("Cantidad reseñas airbnb de Roberto: ", 39),
    @("Precio por noche airbnb de Roberto: ", 26)
)

$row = 4
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
